$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '33.969.33'
$ws.Range("E2").Value = '  -1.94%  '

$ws.Range("D3").Value = '1.790.13'
$ws.Range("E3").Value = '  +0.11%  '

$ws.Range("E4").Value = '  -0.05%  '

$ws.Range("D5").Value = "'222.06"
$ws.Range("E5").Value = '  -0.42%  '

$ws.Range("D6").Value = "'0.550"
$ws.Range("E6").Value = '  -1.04%  '

$ws.Range("E7").Value = '  -0.12%  '

$ws.Range("D8").Value = "'31.54"
$ws.Range("E8").Value = '  -3.40%  '

$ws.Range("E9").Value = '  +1.37%  '

$ws.Range("E10").Value = '  +5.23%  '

$ws.Range("D11").Value = "'0.0921"
$ws.Range("E11").Value = '  -1.64%  '

$ws.Range("D12").Value = '2.046.43'
$ws.Range("E12").Value = '  +0.06%  '

$ws.Range("D13").Value = '1.791.60'
$ws.Range("E13").Value = '  +0.35%  '

$ws.Range("D14").Value = "'10.65"
$ws.Range("E14").Value = '  -4.65%  '

$ws.Range("D15").Value = "'0.629"
$ws.Range("E15").Value = '  -0.43%  '

$ws.Range("D16").Value = '33.935.65'
$ws.Range("E16").Value = '  -1.95%  '

$ws.Range("D17").Value = "'4.23"
$ws.Range("E17").Value = '  -1.92%  '

$ws.Range("D18").Value = "'68.03"
$ws.Range("E18").Value = '  -0.67%  '

$ws.Range("D19").Value = "'245.27"
$ws.Range("E19").Value = '  -3.14%  '

$ws.Range("E20").Value = '  +1.25%  '

$ws.Range("E21").Value = '  +0.11%  '

$ws.Range("D22").Value = "'10.77"
$ws.Range("E22").Value = '  +3.06%  '

$ws.Range("D23").Value = "'4.09"
$ws.Range("E23").Value = '  -3.00%  '

$ws.Range("E24").Value = '  -1.45%  '

$ws.Range("D25").Value = "'158.25"
$ws.Range("E25").Value = '  -0.31%  '

$ws.Range("D26").Value = "'16.40"
$ws.Range("E26").Value = '  +0.33%  '

$ws.Range("D27").Value = "'7.03"
$ws.Range("E27").Value = '  -0.72%  '

$ws.Range("E28").Value = '  -2.02%  '

$ws.Range("E29").Value = '  -0.06%  '

$ws.Range("D30").Value = "'0.0520"
$ws.Range("E30").Value = '  +0.88%  '

$ws.Range("D31").Value = "'1.21"
$ws.Range("E31").Value = '  +1.21%  '

$ws.Range("D32").Value = "'3.70"
$ws.Range("E32").Value = '  -1.50%  '

$ws.Range("D33").Value = "'3.51"
$ws.Range("E33").Value = '  -1.87%  '

$ws.Range("E34").Value = '  -1.53%  '

$ws.Range("D35").Value = '1.409.75'
$ws.Range("E35").Value = '  -2.15%  '

$ws.Range("E36").Value = '  +1.76%  '

$ws.Range("E37").Value = '  -0.01%  '

$ws.Range("D39").Value = "'0.944"
$ws.Range("E39").Value = '  +4.19%  '

$ws.Range("D40").Value = "'79.82"
$ws.Range("E40").Value = '  -3.94%  '

$ws.Range("E41").Value = '  -3.06%  '

$ws.Range("E42").Value = '  -0.49%  '

$ws.Range("D43").Value = "'2.12"
$ws.Range("E43").Value = '  +2.08%  '

$ws.Range("D44").Value = "'5.95"
$ws.Range("E44").Value = '  -0.06%  '

$ws.Range("E45").Value = '  -2.25%  '

$ws.Range("D46").Value = '1.945.56'
$ws.Range("E46").Value = '  -0.01%  '

$ws.Range("E47").Value = '  -0.76%  '

$ws.Range("D48").Value = "'105.57"
$ws.Range("E48").Value = '  +0.46%  '

$ws.Range("D49").Value = "'0.997"
$ws.Range("E49").Value = '  -0.27%  '

$ws.Range("D50").Value = "'11.86"
$ws.Range("E50").Value = '  -1.43%  '

$ws.Range("E51").Value = '  -0.39%  '
